$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.3655573333333333
$ws.Range("N2").Value = 1.096672
$ws.Range("O2").Value = 0.2375364113931583
$ws.Range("P2").Value = 0.2572707166041556
$ws.Range("Q2").Value = 0.1419829556764445
$ws.Range("R2").Value = 1.277846601088
$ws.Range("S2").Value = 0.2375364113931583
$ws.Range("T2").Value = 0.2572707166041556

# Row 3
$ws.Range("M3").Value = 0.6295006666666667
$ws.Range("O3").Value = 0.4090448082825151
$ws.Range("P3").Value = 0.4430278723705731
$ws.Range("S3").Value = 0.4090448082825151
$ws.Range("T3").Value = 0.4430278723705731

# Row 4
$ws.Range("M4").Value = 0.1127876666666667
$ws.Range("N4").Value = 0.338363
$ws.Range("O4").Value = 0.07328857923629238
$ws.Range("P4").Value = 0.07937732656831935
$ws.Range("Q4").Value = 0.04380688011688889
$ws.Range("R4").Value = 0.394261921052
$ws.Range("S4").Value = 0.07328857923629238
$ws.Range("T4").Value = 0.07937732656831935

# Row 5
$ws.Range("M5").Value = 0.3541425
$ws.Range("N5").Value = 0.7082850000000001
$ws.Range("O5").Value = 0.2301191383708208
$ws.Range("P5").Value = 0.1661581489360305
$ws.Range("Q5").Value = 0.13754941919
$ws.Range("R5").Value = 0.8252965151400001
$ws.Range("S5").Value = 0.2301191383708208
$ws.Range("T5").Value = 0.1661581489360305

# Row 6
$ws.Range("M6").Value = 0.07696466666666667
$ws.Range("N6").Value = 0.230894
$ws.Range("O6").Value = 0.05001106271721345
$ws.Range("P6").Value = 0.0541659355209214
$ws.Range("Q6").Value = 0.02989317915288889
$ws.Range("R6").Value = 0.269038612376
$ws.Range("S6").Value = 0.05001106271721345
$ws.Range("T6").Value = 0.0541659355209214
